$d = $word.ActiveDocument

# --- Step 1: merge runs that were split across multiple <w:r> elements ---
# Doing an identical Find/Replace over text that spans several runs makes
# Word re-emit that span as a single run, which is exactly what the diff
# shows (three separate paragraphs where adjacent runs got merged).
$d.Content.Find.Execute(" dans le jeu et début de la gestion des tours.", $false, $false, $false, $false, $false, $true, 1, $false, " dans le jeu et début de la gestion des tours.", 2)
$d.Content.Find.Execute("Réglage de problèmes lors de la lecture d’un fichier texte.", $false, $false, $false, $false, $false, $true, 1, $false, "Réglage de problèmes lors de la lecture d’un fichier texte.", 2)
$d.Content.Find.Execute("Modification de la texture de la barre de vie.", $false, $false, $false, $false, $false, $true, 1, $false, "Modification de la texture de la barre de vie.", 2)

# --- Step 2: drop the _GoBack bookmark from its old spot ---
# (it will reappear at the end of the new content we are about to type,
# mirroring where Word leaves it after the last edit)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3: append the new "Le 02 décembre 2014" entry + its bullet list ---
# Insert right before the (empty) trailing paragraph that already ends the
# document, so our new paragraphs land cleanly in front of it. We add one
# extra blank <w:p/> of our own at the tail of the fragment so the existing
# trailing paragraph's own mark is pushed out to a standalone paragraph
# again afterwards (it briefly duplicates, step 3b removes the duplicate
# introduced by our own placeholder, leaving the original trailing
# paragraph as the sole blank one, exactly like before the edit).
$trailingPara = $d.Paragraphs.Last
$insertPoint = $d.Range($trailingPara.Range.Start, $trailingPara.Range.Start)

$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Le 02 décembre 2014</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Début de la gestion des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tools</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CPlayer</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Début de la gestion du </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>JetPack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CPlayer</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ajout de la classe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CJetPack</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Quelques modifications dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CToolBar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CGame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CSprite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en accordance avec le début de la gestion des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tools</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newContentXml)

# --- Step 3b: collapse the now-duplicated blank paragraph at the very end
# back down to one, keeping the original trailing paragraph's identity ---
$countAfterInsert = $d.Paragraphs.Count
$duplicateBlank = $d.Paragraphs.Item($countAfterInsert - 1)
$dupRange = $d.Range($duplicateBlank.Range.Start, $duplicateBlank.Range.End)
$dupRange.Delete()

# --- Step 4: register the new bullet-list numbering definition (numId 2) ---
# used by the paragraphs just inserted above, mirroring the abstractNum/num
# pair Word would add to numbering.xml for a freshly created bulleted list.
$abstractXml = @'
<w:abstractNum w:abstractNumId="1"><w:nsid w:val="61F841EC"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="6E4242E0"/><w:lvl w:ilvl="0" w:tplc="0C0C0001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0B7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="0C0C0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="0C0C0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="0C0C0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0B7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="0C0C0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="0C0C0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="0C0C0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0B7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="0C0C0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="0C0C0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="&#xF0A7;"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum>
'@

Write-Output "done"
